# Updates cryptos list price (D) and 1h-volume-change (E) columns
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.007.64"
$ws.Range("E2").Value = "  -0.34%  "
$ws.Range("D3").Value = "1.638.19"
$ws.Range("E3").Value = "  -0.57%  "
$ws.Range("E4").Value = "  -0.30%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "214.80"
$ws.Range("E5").Value = "  -1.09%  "
$ws.Range("E6").Value = "  -0.87%  "
$ws.Range("E7").Value = "  -0.21%  "
$ws.Range("E8").Value = "  -1.75%  "
$ws.Range("E9").Value = "  -2.64%  "
$ws.Range("D10").Value = "18.54"
$ws.Range("E10").Value = "  -5.93%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0790"
$ws.Range("E11").Value = "  -0.95%  "
$ws.Range("D12").Value = "1.866.85"
$ws.Range("E12").Value = "  -0.42%  "
$ws.Range("D13").Value = "1.648.79"
$ws.Range("E13").Value = "  +0.83%  "
$ws.Range("E14").Value = "  -1.89%  "
$ws.Range("E15").Value = "  -2.89%  "
$ws.Range("D16").Value = "26.020.39"
$ws.Range("E16").Value = "  -0.54%  "
$ws.Range("E17").Value = "  -2.66%  "
$ws.Range("D18").Value = "61.92"
$ws.Range("E18").Value = "  -2.50%  "
$ws.Range("E19").Value = "  -0.30%  "
$ws.Range("D20").Value = "192.47"
$ws.Range("E20").Value = "  -0.78%  "
$ws.Range("D21").Value = "4.26"
$ws.Range("E21").Value = "  -2.25%  "
$ws.Range("D22").Value = "9.77"
$ws.Range("E22").Value = "  -1.88%  "
$ws.Range("D23").Value = "6.12"
$ws.Range("E23").Value = "  -2.18%  "
$ws.Range("E24").Value = "  +1.90%  "
$ws.Range("D25").Value = "1.79"
$ws.Range("E25").Value = "  -1.61%  "
$ws.Range("D26").Value = "143.78"
$ws.Range("E26").Value = "  -0.59%  "
$ws.Range("E27").Value = "  -0.29%  "
$ws.Range("E28").Value = "  -1.04%  "
$ws.Range("D29").Value = "15.26"
$ws.Range("E29").Value = "  -2.18%  "
$ws.Range("D30").Value = "1.24"
$ws.Range("E30").Value = "  -1.41%  "
$ws.Range("D31").Value = "0.0485"
$ws.Range("E31").Value = "  -2.90%  "
$ws.Range("E32").Value = "  -3.64%  "
$ws.Range("E33").Value = "  -4.73%  "
$ws.Range("E34").Value = "  -2.47%  "
$ws.Range("E35").Value = "  -2.40%  "
$ws.Range("D36").Value = "1.139.92"
$ws.Range("E36").Value = "  +0.31%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.870"
$ws.Range("E37").Value = "  -4.26%  "
$ws.Range("E38").Value = "  -1.05%  "
$ws.Range("D39").Value = "0.522"
$ws.Range("E39").Value = "  -3.68%  "
$ws.Range("E40").Value = "  -1.27%  "
$ws.Range("D41").Value = "98.54"
$ws.Range("E41").Value = "  -1.16%  "
$ws.Range("E42").Value = "  -2.57%  "
$ws.Range("D43").Value = "1.776.16"
$ws.Range("E43").Value = "  -0.35%  "
$ws.Range("D44").Value = "5.24"
$ws.Range("E44").Value = "  -4.89%  "
$ws.Range("E45").Value = "  -1.42%  "
$ws.Range("D46").Value = "55.27"
$ws.Range("E46").Value = "  -2.63%  "
$ws.Range("E47").Value = "  -0.63%  "
$ws.Range("E48").Value = "  +2.42%  "
$ws.Range("E49").Value = "  -0.71%  "
$ws.Range("D50").Value = "7.57"
$ws.Range("E50").Value = "  -2.25%  "
$ws.Range("E51").Value = "  +0.00%  "
